# The diff inserts one new data row right before the current row 49,
# pushing the existing rows 49-87 down to 50-88 (dimension grows from
# A1:R87 to A1:R88). The newly inserted row 49 repeats the same
# "template" columns (A,B,C,E,F,G,H,I,N,O,Q,R) as the row that used to
# sit at 49, but carries its own date/volume/price figures (D,J,K,L,M,P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 49; Excel shifts rows 49:87 down to 50:88.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record.
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = 'Vega Modelo de Temuco'
$ws.Range("C49").Value = 'La Araucanía'
$ws.Range("D49").Value = 44729
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 100112035
$ws.Range("G49").Value = 'Bruselas (repollito)'
$ws.Range("H49").Value = 'Sin especificar'
$ws.Range("I49").Value = 'Primera'
$ws.Range("J49").Value = 45
$ws.Range("K49").Value = 28000
$ws.Range("L49").Value = 28000
$ws.Range("M49").Value = 28000
$ws.Range("N49").Value = '$/malla 10 kilos'
$ws.Range("O49").Value = 'Provincia de Quillota'
$ws.Range("P49").Value = 2800
$ws.Range("Q49").Value = 10
$ws.Range("R49").Value = 'Hortaliza'

# Match the date cell's number format style used by every other date
# cell in column D (style index 2 = "YYYY-MM-DD HH:MM:SS").
$ws.Range("D49").NumberFormat = $ws.Range("D50").NumberFormat
